$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume(1h)/Hora columns to Text number format so
# numeric-looking strings (prices, percentages, hour) are written back
# as text, matching the original inline-string cell contents, instead
# of being auto-coerced to numbers by the Value setter.
$ws.Range("D2:E50").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Every row's "Hora" value moves from 16 to 17.
$ws.Range("G2:G51").Value = "17"

# Per-row Price / Volume(1h) updates.
$ws.Range("D2").Value = "302.67"
$ws.Range("E2").Value = "1.78%"
$ws.Range("D3").Value = "31.84"
$ws.Range("E3").Value = "0.58%"
$ws.Range("D4").Value = "4.990"
$ws.Range("E4").Value = "-1.57%"
$ws.Range("D5").Value = "0.07821"
$ws.Range("E5").Value = "-2.58%"
$ws.Range("D6").Value = "2.120"
$ws.Range("E6").Value = "-19.14%"
$ws.Range("D7").Value = "7.789"
$ws.Range("E7").Value = "-0.07%"
$ws.Range("D8").Value = "3.783"
$ws.Range("E8").Value = "-0.93%"
$ws.Range("D9").Value = "0.9210"
$ws.Range("E9").Value = "-0.36%"
$ws.Range("D10").Value = "0.1747"
$ws.Range("E10").Value = "-0.44%"
$ws.Range("D11").Value = "0.07789"
$ws.Range("E11").Value = "5.58%"
$ws.Range("D12").Value = "0.08893"
$ws.Range("E12").Value = "-0.38%"
$ws.Range("D13").Value = "0.03089"
$ws.Range("E13").Value = "1.00%"
$ws.Range("E14").Value = "0.08%"
$ws.Range("D15").Value = "0.001524"
$ws.Range("E15").Value = "2.08%"
$ws.Range("D16").Value = "0.005805"
$ws.Range("E16").Value = "-4.17%"
$ws.Range("D17").Value = "3.463"
$ws.Range("E17").Value = "-2.43%"
$ws.Range("E18").Value = "0.86%"
$ws.Range("D20").Value = "0.1329"
$ws.Range("E20").Value = "-1.22%"
$ws.Range("D21").Value = "4.148"
$ws.Range("E21").Value = "3.69%"
$ws.Range("D22").Value = "0.1795"
$ws.Range("E22").Value = "9.04%"
$ws.Range("D23").Value = "0.04589"
$ws.Range("E23").Value = "-0.30%"
$ws.Range("D24").Value = "0.001239"
$ws.Range("E24").Value = "0.06%"
$ws.Range("D25").Value = "0.004475"
$ws.Range("E25").Value = "0.96%"
$ws.Range("D26").Value = "0.0001249"
$ws.Range("E26").Value = "4.26%"
$ws.Range("D39").Value = "0.01758"
$ws.Range("E39").Value = "-0.33%"
$ws.Range("D40").Value = "0.04768"
$ws.Range("E40").Value = "6.07%"
$ws.Range("D41").Value = "0.007034"
$ws.Range("E41").Value = "4.10%"
$ws.Range("D42").Value = "0.1368"
$ws.Range("E42").Value = "1.62%"
$ws.Range("D43").Value = "0.002139"
$ws.Range("E43").Value = "-3.08%"
$ws.Range("D44").Value = "0.01074"
$ws.Range("E44").Value = "9.09%"
$ws.Range("D45").Value = "0.00005961"
$ws.Range("E45").Value = "-7.69%"
$ws.Range("E46").Value = "0.16%"
$ws.Range("D47").Value = "0.003551"
$ws.Range("E47").Value = "-59.37%"
$ws.Range("D48").Value = "0.8150"
$ws.Range("E48").Value = "-0.68%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.16%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.16%"

# Reset style back to Normal (drops the explicit Text number format
# we applied above) while the underlying cell values remain text, so
# the saved cells carry no extra style index vs. the original file.
$ws.Range("D2:E50").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
